$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 25
$link = "https://www.biocentury.com/article/656269/illumina-s-buy-of-somalogic-could-be-tipping-point-for-multiomics-deals-report"
$keyword = "BCMA"
$title = "Illumina`u{2019}s buy of Somalogic could be tipping point for multiomics: Deals Report"

$ws.Cells.Item($newRow, 1).Value = $link
$ws.Cells.Item($newRow, 2).Value = $keyword
$ws.Cells.Item($newRow, 3).Value = $title

$ws.Hyperlinks.Add($ws.Cells.Item($newRow, 1), $link) | Out-Null
